$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 7 new data rows (9-15) that replicate the existing "ranking"
# repository entry (rows 2-8), continuing the id sequence 8..14.
# Copy formatting + values from the last existing row (row 8) and then
# patch the id column with the correct sequential value.
$srcRow = $ws.Range("A8:O8")

for ($r = 9; $r -le 15; $r++) {
    $dstRow = $ws.Range("A" + $r + ":O" + $r)
    $srcRow.Copy($dstRow)
}

for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
}
